$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 107
$ws.Cells.Item(107, 8).Value = 1141.1765
$ws.Cells.Item(107, 9).Value = 622
$ws.Cells.Item(107, 11).Value = 622
$ws.Cells.Item(107, 13).Value = 1298

# row 132
$ws.Cells.Item(132, 8).Value = 753.27026
$ws.Cells.Item(132, 9).Value = 636.53125
$ws.Cells.Item(132, 10).Value = 1500.4
$ws.Cells.Item(132, 11).Value = 1909.59375
$ws.Cells.Item(132, 12).Value = 4501.200000000001
$ws.Cells.Item(132, 13).Value = 620.40625
$ws.Cells.Item(132, 14).Value = -9561.200000000001

# row 136
$ws.Cells.Item(136, 8).Value = 70207.8
$ws.Cells.Item(136, 10).Value = 70207.8
$ws.Cells.Item(136, 12).Value = 70207.8
$ws.Cells.Item(136, 14).Value = -80407.8

# row 137
$ws.Cells.Item(137, 8).Value = 2262.4119
$ws.Cells.Item(137, 9).Value = 2055.5454
$ws.Cells.Item(137, 10).Value = 2641.6667
$ws.Cells.Item(137, 11).Value = 6166.6362
$ws.Cells.Item(137, 12).Value = 7925.000100000001
$ws.Cells.Item(137, 13).Value = -3616.6362
$ws.Cells.Item(137, 14).Value = -13025.0001

# row 138
$ws.Cells.Item(138, 8).Value = 1787.871
$ws.Cells.Item(138, 9).Value = 1308.44
$ws.Cells.Item(138, 10).Value = 2111.8108
$ws.Cells.Item(138, 11).Value = 3925.32
$ws.Cells.Item(138, 12).Value = 6335.432400000001
$ws.Cells.Item(138, 13).Value = 1214.68
$ws.Cells.Item(138, 14).Value = -16615.4324

# row 140
$ws.Cells.Item(140, 8).Value = 71142.89999999999
$ws.Cells.Item(140, 10).Value = 71142.89999999999
$ws.Cells.Item(140, 12).Value = 71142.89999999999
$ws.Cells.Item(140, 14).Value = -81502.89999999999

$ws = $wb.Worksheets.Item("ARM")
# row 32
$ws.Cells.Item(32, 8).Value = 5739.9434
$ws.Cells.Item(32, 9).Value = 3710
$ws.Cells.Item(32, 10).Value = 17158.375
$ws.Cells.Item(32, 11).Value = 3710
$ws.Cells.Item(32, 12).Value = 17158.375
$ws.Cells.Item(32, 13).Value = -3423
$ws.Cells.Item(32, 14).Value = -17732.375

# row 44
$ws.Cells.Item(44, 8).Value = 30000
$ws.Cells.Item(44, 10).Value = 30000
$ws.Cells.Item(44, 12).Value = 30000
$ws.Cells.Item(44, 14).Value = -30976

# row 45
$ws.Cells.Item(45, 8).Value = 4738165
$ws.Cells.Item(45, 9).Value = 10000839
$ws.Cells.Item(45, 11).Value = 10000839
$ws.Cells.Item(45, 13).Value = -10000462

# row 61
$ws.Cells.Item(61, 8).Value = 6587.227
$ws.Cells.Item(61, 9).Value = 7901.769
$ws.Cells.Item(61, 11).Value = 7901.769
$ws.Cells.Item(61, 13).Value = -7689.769

# row 122
$ws.Cells.Item(122, 8).Value = 1148.0476
$ws.Cells.Item(122, 9).Value = 849.9286
$ws.Cells.Item(122, 11).Value = 2549.7858
$ws.Cells.Item(122, 13).Value = -99.78579999999965

# row 132
$ws.Cells.Item(132, 8).Value = 1563.6
$ws.Cells.Item(132, 9).Value = 1230.1852
$ws.Cells.Item(132, 11).Value = 3690.5556
$ws.Cells.Item(132, 13).Value = -1160.5556

# row 136
$ws.Cells.Item(136, 8).Value = 6587.227
$ws.Cells.Item(136, 9).Value = 7901.769
$ws.Cells.Item(136, 11).Value = 23705.307
$ws.Cells.Item(136, 13).Value = -21155.307

$ws = $wb.Worksheets.Item("BSM")
# row 20
$ws.Cells.Item(20, 8).Value = 4750
$ws.Cells.Item(20, 9).Value = 4500
$ws.Cells.Item(20, 11).Value = 4500
$ws.Cells.Item(20, 13).Value = -4253

# row 86
$ws.Cells.Item(86, 8).Value = 227145.44
$ws.Cells.Item(86, 9).Value = 8199.75
$ws.Cells.Item(86, 10).Value = 402302
$ws.Cells.Item(86, 11).Value = 8199.75
$ws.Cells.Item(86, 12).Value = 402302
$ws.Cells.Item(86, 13).Value = -7076.75
$ws.Cells.Item(86, 14).Value = -404548

# row 89
$ws.Cells.Item(89, 8).Value = 227145.44
$ws.Cells.Item(89, 9).Value = 8199.75
$ws.Cells.Item(89, 10).Value = 402302
$ws.Cells.Item(89, 11).Value = 40998.75
$ws.Cells.Item(89, 12).Value = 2011510
$ws.Cells.Item(89, 13).Value = -35382.75
$ws.Cells.Item(89, 14).Value = -2022742

# row 134
$ws.Cells.Item(134, 8).Value = 5399.107
$ws.Cells.Item(134, 10).Value = 2953.5
$ws.Cells.Item(134, 12).Value = 8860.5
$ws.Cells.Item(134, 14).Value = -13930.5

$ws = $wb.Worksheets.Item("CRP")
# row 31
$ws.Cells.Item(31, 8).Value = 2316.611
$ws.Cells.Item(31, 9).Value = 2199.75
$ws.Cells.Item(31, 11).Value = 2199.75
$ws.Cells.Item(31, 13).Value = -1904.75

# row 34
$ws.Cells.Item(34, 8).Value = 2316.611
$ws.Cells.Item(34, 9).Value = 2199.75
$ws.Cells.Item(34, 11).Value = 2199.75
$ws.Cells.Item(34, 13).Value = -1997.75

# row 58
$ws.Cells.Item(58, 8).Value = 4833519
$ws.Cells.Item(58, 9).Value = 10870368
$ws.Cells.Item(58, 11).Value = 10870368
$ws.Cells.Item(58, 13).Value = -10870165

# row 132
$ws.Cells.Item(132, 8).Value = 2276.8845
$ws.Cells.Item(132, 9).Value = 1431.1052
$ws.Cells.Item(132, 11).Value = 4293.3156
$ws.Cells.Item(132, 13).Value = -1763.3156

# row 134
$ws.Cells.Item(134, 8).Value = 3233.3076
$ws.Cells.Item(134, 9).Value = 2934.111
$ws.Cells.Item(134, 11).Value = 8802.332999999999
$ws.Cells.Item(134, 13).Value = -6267.332999999999

# row 136
$ws.Cells.Item(136, 8).Value = 4833519
$ws.Cells.Item(136, 9).Value = 10870368
$ws.Cells.Item(136, 11).Value = 32611104
$ws.Cells.Item(136, 13).Value = -32608554

$ws = $wb.Worksheets.Item("CUL")
# row 12
$ws.Cells.Item(12, 8).Value = 90.30768999999999
$ws.Cells.Item(12, 9).Value = 66.59999999999999
$ws.Cells.Item(12, 10).Value = 105.125
$ws.Cells.Item(12, 11).Value = 199.8
$ws.Cells.Item(12, 12).Value = 315.375
$ws.Cells.Item(12, 13).Value = -26.79999999999998
$ws.Cells.Item(12, 14).Value = -661.375

# row 131
$ws.Cells.Item(131, 8).Value = 15284.236
$ws.Cells.Item(131, 10).Value = 15832.132
$ws.Cells.Item(131, 12).Value = 47496.396
$ws.Cells.Item(131, 14).Value = -57576.396

# row 136
$ws.Cells.Item(136, 8).Value = 1415.9
$ws.Cells.Item(136, 9).Value = 1415.9
$ws.Cells.Item(136, 11).Value = 4247.700000000001
$ws.Cells.Item(136, 13).Value = 852.2999999999993

# row 140
$ws.Cells.Item(140, 8).Value = 3284.5715
$ws.Cells.Item(140, 9).Value = 1212.2858
$ws.Cells.Item(140, 10).Value = 5356.857
$ws.Cells.Item(140, 11).Value = 3636.8574
$ws.Cells.Item(140, 12).Value = 16070.571
$ws.Cells.Item(140, 13).Value = 1543.1426
$ws.Cells.Item(140, 14).Value = -26430.571

$ws = $wb.Worksheets.Item("GSM")
# row 122
$ws.Cells.Item(122, 8).Value = 1630.0625
$ws.Cells.Item(122, 9).Value = 1488.5454
$ws.Cells.Item(122, 11).Value = 4465.6362
$ws.Cells.Item(122, 13).Value = -2015.6362

# row 126
$ws.Cells.Item(126, 8).Value = 2695780.2
$ws.Cells.Item(126, 10).Value = 102278.8
$ws.Cells.Item(126, 12).Value = 306836.4
$ws.Cells.Item(126, 14).Value = -311776.4

# row 136
$ws.Cells.Item(136, 8).Value = 10158.6
$ws.Cells.Item(136, 10).Value = 10158.6
$ws.Cells.Item(136, 12).Value = 30475.8
$ws.Cells.Item(136, 14).Value = -35575.8

$ws = $wb.Worksheets.Item("LTW")
# row 16
$ws.Cells.Item(16, 8).Value = 8310.375
$ws.Cells.Item(16, 9).Value = 9354.714
$ws.Cells.Item(16, 10).Value = 1000
$ws.Cells.Item(16, 11).Value = 9354.714
$ws.Cells.Item(16, 12).Value = 1000
$ws.Cells.Item(16, 13).Value = -9184.714
$ws.Cells.Item(16, 14).Value = -1340

# row 26
$ws.Cells.Item(26, 8).Value = 8420
$ws.Cells.Item(26, 9).Value = 0
$ws.Cells.Item(26, 11).Value = 0
$ws.Cells.Item(26, 13).ClearContents()

# row 40
$ws.Cells.Item(40, 8).Value = 11017.059
$ws.Cells.Item(40, 9).Value = 11283
$ws.Cells.Item(40, 11).Value = 11283
$ws.Cells.Item(40, 13).Value = -11147

# row 61
$ws.Cells.Item(61, 8).Value = 2314.0527
$ws.Cells.Item(61, 9).Value = 2097.25
$ws.Cells.Item(61, 11).Value = 2097.25
$ws.Cells.Item(61, 13).Value = -1895.25

# row 113
$ws.Cells.Item(113, 8).Value = 2314.0527
$ws.Cells.Item(113, 9).Value = 2097.25
$ws.Cells.Item(113, 11).Value = 2097.25
$ws.Cells.Item(113, 13).Value = 72.75

# row 122
$ws.Cells.Item(122, 8).Value = 6464.5884
$ws.Cells.Item(122, 9).Value = 5657.75
$ws.Cells.Item(122, 11).Value = 16973.25
$ws.Cells.Item(122, 13).Value = -14523.25

# row 132
$ws.Cells.Item(132, 8).Value = 1845.1777
$ws.Cells.Item(132, 9).Value = 1265
$ws.Cells.Item(132, 11).Value = 3795
$ws.Cells.Item(132, 13).Value = -1265

# row 135
$ws.Cells.Item(135, 8).Value = 35812.5
$ws.Cells.Item(135, 10).Value = 35812.5
$ws.Cells.Item(135, 12).Value = 35812.5
$ws.Cells.Item(135, 14).Value = -45952.5

# row 136
$ws.Cells.Item(136, 8).Value = 4558.3
$ws.Cells.Item(136, 9).Value = 3262.3333
$ws.Cells.Item(136, 11).Value = 9786.999899999999
$ws.Cells.Item(136, 13).Value = -7236.999899999999

$ws = $wb.Worksheets.Item("WVR")
# row 46
$ws.Cells.Item(46, 8).Value = 0
$ws.Cells.Item(46, 10).Value = 0
$ws.Cells.Item(46, 12).Value = 0
$ws.Cells.Item(46, 14).ClearContents()

# row 81
$ws.Cells.Item(81, 8).Value = 1633.8334
$ws.Cells.Item(81, 10).Value = 402.5
$ws.Cells.Item(81, 12).Value = 805
$ws.Cells.Item(81, 14).Value = -2927

# row 84
$ws.Cells.Item(84, 8).Value = 1633.8334
$ws.Cells.Item(84, 10).Value = 402.5
$ws.Cells.Item(84, 12).Value = 4025
$ws.Cells.Item(84, 14).Value = -14633

# row 132
$ws.Cells.Item(132, 8).Value = 1039
$ws.Cells.Item(132, 9).Value = 836.6326
$ws.Cells.Item(132, 11).Value = 2509.8978
$ws.Cells.Item(132, 13).Value = 20.10219999999981

# row 134
$ws.Cells.Item(134, 8).Value = 0
$ws.Cells.Item(134, 10).Value = 0
$ws.Cells.Item(134, 12).Value = 0
$ws.Cells.Item(134, 14).ClearContents()
